$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("G2").Value = 40.7349555
$ws.Range("H2").Value = 81.469911
$ws.Range("I2").Value = 0.05567871843833241
$ws.Range("J2").Value = 0.03826666865920979
$ws.Range("K2").Value = 2
$ws.Range("M2").Value = 0.4274105
$ws.Range("N2").Value = 0.854821
$ws.Range("O2").Value = 0.4422185568930271
$ws.Range("P2").Value = 0.3984644485961792
$ws.Range("Q2").Value = 17.41054769773275
$ws.Range("R2").Value = 69.642190790931
$ws.Range("S2").Value = 0.02462216251745254
$ws.Range("T2").Value = 0.01524790702690472

# Row 3
$ws.Range("D3").Value = "Neutro"
$ws.Range("E3").Value = 2
$ws.Range("G3").Value = 40.7349555
$ws.Range("H3").Value = 81.469911
$ws.Range("I3").Value = 0.05567871843833241
$ws.Range("J3").Value = 0.03826666865920979
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.2122596666666667
$ws.Range("N3").Value = 0.636779
$ws.Range("O3").Value = 0.2196136114576696
$ws.Range("P3").Value = 0.2968268129966699
$ws.Range("Q3").Value = 8.6463880761115
$ws.Range("R3").Value = 51.878328456669
$ws.Range("S3").Value = 0.01222780443757692
$ws.Range("T3").Value = 0.01135857330211279

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("G4").Value = 40.7349555
$ws.Range("H4").Value = 81.469911
$ws.Range("I4").Value = 0.05567871843833241
$ws.Range("J4").Value = 0.03826666865920979
$ws.Range("K4").Value = 2
$ws.Range("M4").Value = 0.326844
$ws.Range("N4").Value = 0.653688
$ws.Range("O4").Value = 0.3381678316493033
$ws.Range("P4").Value = 0.304708738407151
$ws.Range("Q4").Value = 13.313975795442
$ws.Range("R4").Value = 53.255903181768
$ws.Range("S4").Value = 0.01882875148330295
$ws.Range("T4").Value = 0.01166018833019228

# Row 5
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 350.3919066666667
$ws.Range("H5").Value = 1051.17572
$ws.Range("I5").Value = 0.4789344206933965
$ws.Range("J5").Value = 0.4937404802104949
$ws.Range("K5").Value = 2
$ws.Range("M5").Value = 0.4274105
$ws.Range("N5").Value = 0.854821
$ws.Range("O5").Value = 0.4422185568930271
$ws.Range("P5").Value = 0.3984644485961792
$ws.Range("Q5").Value = 149.7611800243533
$ws.Range("R5").Value = 898.56708014612
$ws.Range("S5").Value = 0.2117936883654317
$ws.Range("T5").Value = 0.1967380281966876

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("D6").Value = "Neutro"
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 350.3919066666667
$ws.Range("H6").Value = 1051.17572
$ws.Range("I6").Value = 0.4789344206933965
$ws.Range("J6").Value = 0.4937404802104949
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.2122596666666667
$ws.Range("N6").Value = 0.636779
$ws.Range("O6").Value = 0.2196136114576696
$ws.Range("P6").Value = 0.2968268129966699
$ws.Range("Q6").Value = 74.37406931176444
$ws.Range("R6").Value = 669.36662380588
$ws.Range("S6").Value = 0.1051805177798637
$ws.Range("T6").Value = 0.1465554131883266

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 350.3919066666667
$ws.Range("H7").Value = 1051.17572
$ws.Range("I7").Value = 0.4789344206933965
$ws.Range("J7").Value = 0.4937404802104949
$ws.Range("K7").Value = 2
$ws.Range("M7").Value = 0.326844
$ws.Range("N7").Value = 0.653688
$ws.Range("O7").Value = 0.3381678316493033
$ws.Range("P7").Value = 0.304708738407151
$ws.Range("Q7").Value = 114.52349234256
$ws.Range("R7").Value = 687.14095405536
$ws.Range("S7").Value = 0.1619602145481011
$ws.Range("T7").Value = 0.1504470388254808

# Row 8
$ws.Range("A8").Value = "M1"
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 243.8287033333334
$ws.Range("H8").Value = 731.48611
$ws.Range("I8").Value = 0.3332781281688242
$ws.Range("J8").Value = 0.3435812836494235
$ws.Range("K8").Value = 2
$ws.Range("M8").Value = 0.4274105
$ws.Range("N8").Value = 0.854821
$ws.Range("O8").Value = 0.4422185568930271
$ws.Range("P8").Value = 0.3984644485961792
$ws.Range("Q8").Value = 104.2149480060517
$ws.Range("R8").Value = 625.28968803631
$ws.Range("S8").Value = 0.1473817728828268
$ws.Range("T8").Value = 0.136904926737335

# Row 9
$ws.Range("A9").Value = "M1"
$ws.Range("D9").Value = "Neutro"
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 243.8287033333334
$ws.Range("H9").Value = 731.48611
$ws.Range("I9").Value = 0.3332781281688242
$ws.Range("J9").Value = 0.3435812836494235
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.2122596666666667
$ws.Range("N9").Value = 0.636779
$ws.Range("O9").Value = 0.2196136114576696
$ws.Range("P9").Value = 0.2968268129966699
$ws.Range("Q9").Value = 51.75499929329889
$ws.Range("R9").Value = 465.79499363969
$ws.Range("S9").Value = 0.07319241334700757
$ws.Range("T9").Value = 0.1019841374309632

# Row 10
$ws.Range("A10").Value = "M1"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 243.8287033333334
$ws.Range("H10").Value = 731.48611
$ws.Range("I10").Value = 0.3332781281688242
$ws.Range("J10").Value = 0.3435812836494235
$ws.Range("K10").Value = 2
$ws.Range("M10").Value = 0.326844
$ws.Range("N10").Value = 0.653688
$ws.Range("O10").Value = 0.3381678316493033
$ws.Range("P10").Value = 0.304708738407151
$ws.Range("Q10").Value = 79.69394871228
$ws.Range("R10").Value = 478.1636922736801
$ws.Range("S10").Value = 0.1127039419389899
$ws.Range("T10").Value = 0.1046922194811253

# Row 11
$ws.Range("A11").Value = "M2"
$ws.Range("D11").Value = "ECs"
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 48.737294
$ws.Range("H11").Value = 146.211882
$ws.Range("I11").Value = 0.06661674320651284
$ws.Range("J11").Value = 0.06867617226847689
$ws.Range("K11").Value = 2
$ws.Range("M11").Value = 0.4274105
$ws.Range("N11").Value = 0.854821
$ws.Range("O11").Value = 0.4422185568930271
$ws.Range("P11").Value = 0.3984644485961792
$ws.Range("Q11").Value = 20.830831197187
$ws.Range("R11").Value = 124.984987183122
$ws.Range("S11").Value = 0.02945916004569748
$ws.Range("T11").Value = 0.02736501311465486

# Row 12
$ws.Range("A12").Value = "M2"
$ws.Range("D12").Value = "Neutro"
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 48.737294
$ws.Range("H12").Value = 146.211882
$ws.Range("I12").Value = 0.06661674320651284
$ws.Range("J12").Value = 0.06867617226847689
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.2122596666666667
$ws.Range("N12").Value = 0.636779
$ws.Range("O12").Value = 0.2196136114576696
$ws.Range("P12").Value = 0.2968268129966699
$ws.Range("Q12").Value = 10.34496177867533
$ws.Range("R12").Value = 93.104656008078
$ws.Range("S12").Value = 0.01462994355913046
$ws.Range("T12").Value = 0.02038492934326228

# Row 13
$ws.Range("A13").Value = "M2"
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 48.737294
$ws.Range("H13").Value = 146.211882
$ws.Range("I13").Value = 0.06661674320651284
$ws.Range("J13").Value = 0.06867617226847689
$ws.Range("K13").Value = 2
$ws.Range("M13").Value = 0.326844
$ws.Range("N13").Value = 0.653688
$ws.Range("O13").Value = 0.3381678316493033
$ws.Range("P13").Value = 0.304708738407151
$ws.Range("Q13").Value = 15.929492120136
$ws.Range("R13").Value = 95.57695272081601
$ws.Range("S13").Value = 0.02252763960168491
$ws.Range("T13").Value = 0.02092622981055976

# Row 14 (new)
$ws.Range("A14").Value = "Neutro"
$ws.Range("B14").Value = "Fn1"
$ws.Range("C14").Value = "Col13a1"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 22.832077
$ws.Range("H14").Value = 68.496231
$ws.Range("I14").Value = 0.03120810544755168
$ws.Range("J14").Value = 0.03217289111905
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 0.4274105
$ws.Range("N14").Value = 0.854821
$ws.Range("O14").Value = 0.4422185568930271
$ws.Range("P14").Value = 0.3984644485961792
$ws.Range("Q14").Value = 9.7586694466085
$ws.Range("R14").Value = 58.552016679651
$ws.Range("S14").Value = 0.01380080335438172
$ws.Range("T14").Value = 0.01281975331949717

# Row 15 (new)
$ws.Range("A15").Value = "Neutro"
$ws.Range("B15").Value = "Fn1"
$ws.Range("C15").Value = "Col13a1"
$ws.Range("D15").Value = "Neutro"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 22.832077
$ws.Range("H15").Value = 68.496231
$ws.Range("I15").Value = 0.03120810544755168
$ws.Range("J15").Value = 0.03217289111905
$ws.Range("K15").Value = 2
$ws.Range("L15").Value = 0.6666666666666666
$ws.Range("M15").Value = 0.2122596666666667
$ws.Range("N15").Value = 0.636779
$ws.Range("O15").Value = 0.2196136114576696
$ws.Range("P15").Value = 0.2968268129966699
$ws.Range("Q15").Value = 4.846329053327666
$ws.Range("R15").Value = 43.616961479949
$ws.Range("S15").Value = 0.006853724744088597
$ws.Range("T15").Value = 0.009549776735756476

# Row 16 (new)
$ws.Range("A16").Value = "Neutro"
$ws.Range("B16").Value = "Fn1"
$ws.Range("C16").Value = "Col13a1"
$ws.Range("D16").Value = "sCs"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 22.832077
$ws.Range("H16").Value = 68.496231
$ws.Range("I16").Value = 0.03120810544755168
$ws.Range("J16").Value = 0.03217289111905
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.326844
$ws.Range("N16").Value = 0.653688
$ws.Range("O16").Value = 0.3381678316493033
$ws.Range("P16").Value = 0.304708738407151
$ws.Range("Q16").Value = 7.462527374988
$ws.Range("R16").Value = 44.775164249928
$ws.Range("S16").Value = 0.01055357734908136
$ws.Range("T16").Value = 0.009803361063796357

# Row 17 (new)
$ws.Range("A17").Value = "sCs"
$ws.Range("B17").Value = "Fn1"
$ws.Range("C17").Value = "Col13a1"
$ws.Range("D17").Value = "ECs"
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 25.082339
$ws.Range("H17").Value = 50.164678
$ws.Range("I17").Value = 0.03428388404538221
$ws.Range("J17").Value = 0.02356250409334498
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.4274105
$ws.Range("N17").Value = 0.854821
$ws.Range("O17").Value = 0.4422185568930271
$ws.Range("P17").Value = 0.3984644485961792
$ws.Range("Q17").Value = 10.7204550531595
$ws.Range("R17").Value = 42.881820212638
$ws.Range("S17").Value = 0.0151609697272368
$ws.Range("T17").Value = 0.00938882020109992

# Row 18 (new)
$ws.Range("A18").Value = "sCs"
$ws.Range("B18").Value = "Fn1"
$ws.Range("C18").Value = "Col13a1"
$ws.Range("D18").Value = "Neutro"
$ws.Range("E18").Value = 2
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 25.082339
$ws.Range("H18").Value = 50.164678
$ws.Range("I18").Value = 0.03428388404538221
$ws.Range("J18").Value = 0.02356250409334498
$ws.Range("K18").Value = 2
$ws.Range("L18").Value = 0.6666666666666666
$ws.Range("M18").Value = 0.2122596666666667
$ws.Range("N18").Value = 0.636779
$ws.Range("O18").Value = 0.2196136114576696
$ws.Range("P18").Value = 0.2968268129966699
$ws.Range("Q18").Value = 5.323968915360333
$ws.Range("R18").Value = 31.943813492162
$ws.Range("S18").Value = 0.007529207590002368
$ws.Range("T18").Value = 0.006993982996248579

# Row 19 (new)
$ws.Range("A19").Value = "sCs"
$ws.Range("B19").Value = "Fn1"
$ws.Range("C19").Value = "Col13a1"
$ws.Range("D19").Value = "sCs"
$ws.Range("E19").Value = 2
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 25.082339
$ws.Range("H19").Value = 50.164678
$ws.Range("I19").Value = 0.03428388404538221
$ws.Range("J19").Value = 0.02356250409334498
$ws.Range("K19").Value = 2
$ws.Range("L19").Value = 1
$ws.Range("M19").Value = 0.326844
$ws.Range("N19").Value = 0.653688
$ws.Range("O19").Value = 0.3381678316493033
$ws.Range("P19").Value = 0.304708738407151
$ws.Range("Q19").Value = 8.198012008116
$ws.Range("R19").Value = 32.792048032464
$ws.Range("S19").Value = 0.01159370672814305
$ws.Range("T19").Value = 0.007179700895996478
